# ---------------------------------------------------------------------------
# Applies the "9 - noyabr 2023" revision to Constructors.docx
# ---------------------------------------------------------------------------
$d = $word.ActiveDocument
$RSQUO = [char]0x2019   # the curly apostrophe (U+2019) used throughout the doc

# ---------------------------------------------------------------------------
# 1) Add a first-line (paragraph) indent of 720 twips (36 pt / 0.5")
#    to the "Masalan, ..." paragraph and the "Pastdagi ..." paragraph.
# ---------------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("Masalan, pastdagi misolda")
if ($found) {
    $rng.ParagraphFormat.FirstLineIndent = 36
}

$rng = $d.Content
$found = $rng.Find.Execute("Pastdagi misolda 8-")
if ($found) {
    $rng.ParagraphFormat.FirstLineIndent = 36
}

# ---------------------------------------------------------------------------
# 2) "... 2 la holatda hambizda c.ni recursive ..."
#       -> "... 2 la holatda ham bizda c.ni recursive ..."
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("hambizda", $false, $false, $false, $false, $false, `
    $true, 1, $false, "ham bizda", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3) Insert a new explanatory sentence right after "... runtime error chiqadi"
#    and before the trailing colon, with "StackOverflowError" in bold.
# ---------------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("chiqadi:")
if ($found) {
    $insPos = $rng.End - 1   # just before the ':'

    $newSentence = ". Errorni nomi StackOverflowError bo" + $RSQUO + "lib, 8- va 12-qatorlarda constructorlarimiz o" + $RSQUO + "zini o" + $RSQUO + "zi chaqiryapti. Demak recursive constructor calling bo" + $RSQUO + "lyapti"

    $ins = $d.Range($insPos, $insPos)
    $ins.InsertAfter($newSentence)

    $boldStart = $insPos + ". Errorni nomi ".Length
    $boldEnd = $boldStart + "StackOverflowError".Length
    $boldRange = $d.Range($boldStart, $boldEnd)
    $boldRange.Font.Bold = 1
}

# ---------------------------------------------------------------------------
# 4) Merge the stray empty paragraph that sits right before the
#    "Faqat access modifierlari ..." paragraph into that paragraph, and move
#    the "_GoBack" bookmark there (Word always keeps _GoBack at the location
#    of the most-recent edit).
# ---------------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("Faqat access modifierlari")
if ($found) {
    $faqatPara = $rng.Paragraphs(1)
    $emptyPara = $faqatPara.Previous()
    if ($emptyPara.Range.Text.Trim().Length -eq 0) {
        $emptyPara.Range.Delete()
    }
}

if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$rng = $d.Content
$found = $rng.Find.Execute("Faqat access modifierlari")
if ($found) {
    $faqatPara = $rng.Paragraphs(1)
    $startPos = $faqatPara.Range.Start
    $bmRange = $d.Range($startPos, $startPos)
    $d.Bookmarks.Add("_GoBack", $bmRange)
}

Write-Output "edit complete"
